$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$valsB = @(1.99233694257714,1.829878811781839,1.731337945964469,1.691480612036287,1.684880233897331,1.730799212667876,1.936067247902145,2.348448614481356,2.657823980458943,2.800044957569412,2.854120536436483,2.842464528390167,2.804489351797997,2.781257239095851,2.648559896939332,2.567538907740015,2.521077139257045,2.505369734598844,2.576149263203035,2.815637572992784,2.973439848752264,2.889098373414811,2.572256153549688,2.235793697374675)
for ($i = 0; $i -lt $valsB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $valsB[$i]
}

$valsC = @(0.4809582057064858,0.4314638139891258,0.4013202195111489,0.3890966046190556,0.3870704645300407,0.4011551263416777,0.4638403694764293,0.58880545503348,0.6819905964242707,0.7247103127510854,0.740936749683442,0.7374398656335757,0.7260442693781783,0.7190706335633195,0.6792055637012595,0.6548354285528148,0.6408492440175451,0.6361190029011823,0.6574264611274998,0.7293900747653197,0.7767114245382345,0.7514279998820257,0.6562549790463663,0.5547674565809757)
for ($i = 0; $i -lt $valsC.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $valsC[$i]
}

$valsE = @(0.2795298445871666,0.275144287329077,0.2725784231280883,0.2715647836264239,0.2713984024914069,0.2725646232931567,0.2779914051301233,0.2896382577048939,0.2988068812012514,0.3031106958976508,0.3047595318800447,0.30440357787122,0.3032459645833967,0.3025393758466137,0.2985282883896048,0.2961016376699135,0.2947184098172997,0.2942522239522276,0.2963586636103415,0.3035854663234048,0.3084197638287165,0.3058294506111991,0.296242425157125,0.286380001225794)
for ($i = 0; $i -lt $valsE.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $valsE[$i]
}

$valsF = @(3.228908059506921,3.154345204806646,3.110727570793983,3.093492909977272,3.090663589408479,3.110492957856138,3.202746932335089,3.401040476099126,3.557639581569589,3.631320891495307,3.65957897659888,3.653477149970826,3.633638524835476,3.621533386100424,3.552873946910154,3.511383123767814,3.487748385489624,3.479785361646549,3.515776083890671,3.639455890792817,3.722368654240626,3.677924411094551,3.513789345220715,3.34550353496536)
for ($i = 0; $i -lt $valsF.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $valsF[$i]
}

$valsG = @(0.002478193465876712,0.002484719659774736,0.002488930196346592,0.002490697376749798,0.002490993922837947,0.002488953820908877,0.002480401596202137,0.002465235710152954,0.002455058982524155,0.002450636257392835,0.002448991000044658,0.002449344025288718,0.002450500310216704,0.002451212408538419,0.002455352161104753,0.002457944573806991,0.002459455129751117,0.002459969927957265,0.002457666593593692,0.00245015988138511,0.002445425861599523,0.002447936818192274,0.00245779220568676,0.002469167981714375)
for ($i = 0; $i -lt $valsG.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $valsG[$i]
}

$valsJ = @(0.1368315091867061,0.1332497966659147,0.1311567561595908,0.1303301818044886,0.1301945119387398,0.131145502410476,0.1355743003679422,0.1451171817035544,0.1526743419383365,0.1562359427672391,0.1576028386695754,0.1573076376522522,0.1563480311959182,0.1557626257841775,0.1524441076128653,0.150440291205868,0.149299378852362,0.1489150711666838,0.1506523944275955,0.1566293939156083,0.160641949352879,0.1584905208915615,0.1505564679903273,0.142441374628369)
for ($i = 0; $i -lt $valsJ.Length; $i++) {
    $ws.Cells.Item($i + 2, 10).Value = $valsJ[$i]
}

$valsN = @(1.427685832647192,1.443477760151858,1.453797225831146,1.458158260191652,1.458891787618107,1.45385541078133,1.433000992037719,1.3970915545381,1.373805439223304,1.363896867903406,1.360244246768701,1.361026458914779,1.363594363135704,1.365180275310415,1.374466855176962,1.380339955996789,1.383782358831823,1.384958915595163,1.379708086268657,1.36283739688399,1.352392155261782,1.35791346419866,1.379993549459364,1.406266399016509)
for ($i = 0; $i -lt $valsN.Length; $i++) {
    $ws.Cells.Item($i + 2, 14).Value = $valsN[$i]
}
